# preliminary commit of 3.30.11.00
# Update the "date stamp" cell at top of the change-log sheet and append
# the new rows documenting versions 3.30.10.02 and 3.30.11.00.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the version/date stamp shown at top of sheet (F1)
$ws.Range("F1").Value = "2018-02-27 for 3.30.11.00"

# --- New change-log rows ---

# Row 49: 3.30.10.02 / fix / various
$ws.Range("A49").Value = 43133
$ws.Range("B49").Value = "3.30.10.02"
$ws.Range("C49").Value = "fix"
$ws.Range("D49").Value = "various"
$ws.Range("F49").Value = "fix problem with super year in generalized size comp"
$ws.Range("G49").Value = "No"

# Row 50: 3.30.10.02 / misc / various
$ws.Range("A50").Value = 43133
$ws.Range("B50").Value = "3.30.10.02"
$ws.Range("C50").Value = "misc"
$ws.Range("D50").Value = "various"
$ws.Range("F50").Value = "break SS_write.tpl into SS_write, SS_write_report and SS_write_ssnew"
$ws.Range("G50").Value = "No"

# Row 51: 3.30.10.02 / misc / various
$ws.Range("A51").Value = 43133
$ws.Range("B51").Value = "3.30.10.02"
$ws.Range("C51").Value = "misc"
$ws.Range("D51").Value = "various"
$ws.Range("F51").Value = " clean-up the cout's at end of run"
$ws.Range("G51").Value = "No"

# Row 52: 3.30.11.00 / fix / growth
$ws.Range("A52").Value = 43158
$ws.Range("B52").Value = "3.30.11.00"
$ws.Range("C52").Value = "fix"
$ws.Range("D52").Value = "growth"
$ws.Range("F52").Value = " for Richards growth, disable trap on fish shrinkage due to code interaction"
$ws.Range("G52").Value = "No"

# Row 53: 3.30.11.00 / revise / output
$ws.Range("A53").Value = 43158
$ws.Range("B53").Value = "3.30.11.00"
$ws.Range("C53").Value = "revise"
$ws.Range("D53").Value = "output"
$ws.Range("F53").Value = "revise format of ss_summary"
$ws.Range("G53").Value = "No"

# Row 54: 3.30.11.00 / fix / growth
$ws.Range("A54").Value = 43158
$ws.Range("B54").Value = "3.30.11.00"
$ws.Range("C54").Value = "fix"
$ws.Range("D54").Value = "growth"
$ws.Range("F54").Value = "add totbio, smrybio and totcal catch to end of ss_summary, but without se"
$ws.Range("G54").Value = "No"

# Move the active selection to reflect where editing left off
$ws.Range("F54").Select() | Out-Null
